$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "10÷2=5, 0"
$t.Cell(1,2).Range.Text = "60÷4=15, 0"
$t.Cell(1,3).Range.Text = "92÷5=18, 2"
$t.Cell(1,4).Range.Text = "60÷4=15, 0"
$t.Cell(1,5).Range.Text = "92÷4=23, 0"

$t.Cell(5,1).Range.Text = "77÷4=19, 1"
$t.Cell(5,2).Range.Text = "62÷3=20, 2"
$t.Cell(5,3).Range.Text = "33÷6=5, 3"
$t.Cell(5,4).Range.Text = "69÷4=17, 1"
$t.Cell(5,5).Range.Text = "12÷3=4, 0"

$t.Cell(9,1).Range.Text = "94÷8=11, 6"
$t.Cell(9,2).Range.Text = "38÷3=12, 2"
$t.Cell(9,3).Range.Text = "49÷7=7, 0"
$t.Cell(9,4).Range.Text = "29÷8=3, 5"
$t.Cell(9,5).Range.Text = "81÷3=27, 0"

$t.Cell(13,1).Range.Text = "68÷7=9, 5"
$t.Cell(13,2).Range.Text = "42÷9=4, 6"
$t.Cell(13,3).Range.Text = "90÷8=11, 2"
$t.Cell(13,4).Range.Text = "50÷6=8, 2"
$t.Cell(13,5).Range.Text = "92÷9=10, 2"

$t.Cell(17,1).Range.Text = "81÷8=10, 1"
$t.Cell(17,2).Range.Text = "53÷7=7, 4"
$t.Cell(17,3).Range.Text = "37÷2=18, 1"
$t.Cell(17,4).Range.Text = "12÷7=1, 5"
$t.Cell(17,5).Range.Text = "58÷8=7, 2"
